# Added ability to remove food from ration.
# - Corrects vA (vitamin A) values for the first two foods (H2, H3)
# - Adds three new food rows (Апельсиновый сок / Сыр Рикотта 7-9% / Сыр Моцарелла)
# - Tidies up the sheet view / column widths / page setup left behind by the edit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing values (vA * 10000 column) ------------------------------
$ws.Range("H2").Value = 210000
$ws.Range("H3").Value = 280000

# --- New rows --------------------------------------------------------------
$newFoods = @{
    4 = @{
        "A" = 3; "B" = "Апельсиновый сок"; "C" = 45; "D" = 7000; "E" = 0; "F" = 104000;
        "G" = 11000; "H" = 430000; "I" = 900000; "J" = 300000; "K" = 400; "L" = 190;
        "M" = 400000; "N" = 0; "O" = 300000; "P" = 0; "Q" = 50000; "R" = 0; "S" = 40;
        "T" = 1000; "U" = 200; "V" = 11000; "W" = 0; "X" = 200000; "Y" = 11000; "Z" = 44;
        "AA" = 17000; "AB" = 0; "AC" = 50; "AD" = 100; "AE" = 1000; "AF" = 14
    }
    5 = @{
        "A" = 4; "B" = "Сыр Рикотта 7-9%"; "C" = 140; "D" = 114000; "E" = 79100; "F" = 51000;
        "G" = 70000; "H" = 1250000; "I" = 210000; "J" = 1850000; "K" = 78; "L" = 242;
        "M" = 200000; "N" = 0; "O" = 230000; "P" = 2900; "Q" = 0; "R" = 1000; "S" = 70;
        "T" = 7000; "U" = 440; "V" = 272000; "W" = 0; "X" = 125000; "Y" = 15000; "Z" = 34;
        "AA" = 183000; "AB" = 0; "AC" = 1340; "AD" = 16700; "AE" = 99000; "AF" = 10
    }
    6 = @{
        "A" = 5; "B" = "Сыр Моцарелла"; "C" = 301; "D" = 246000; "E" = 197200; "F" = 64000;
        "G" = 166000; "H" = 2540000; "I" = 590000; "J" = 3480000; "K" = 144; "L" = 415;
        "M" = 1110000; "N" = 0; "O" = 270000; "P" = 18200; "Q" = 0; "R" = 4000; "S" = 430;
        "T" = 13000; "U" = 230; "V" = 716000; "W" = 0; "X" = 131000; "Y" = 29000; "Z" = 34;
        "AA" = 537000; "AB" = 0; "AC" = 3610; "AD" = 26800; "AE" = 682000; "AF" = 41
    }
}

foreach ($rowNum in $newFoods.Keys) {
    $rowValues = $newFoods[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}

# --- Column width tweaks (closest achievable widths in this engine) --------
$ws.Columns("H").ColumnWidth = 11.666666666666666
$ws.Columns("V").ColumnWidth = 10.666666666666666
$ws.Columns("AE").ColumnWidth = 10.666666666666666

# --- Page setup --------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection / view state --------------------------------------------------
$ws.Range("AB6").Select() | Out-Null
